# Add a new "Correction" column (N) to the Card18 sheet, matching the
# pattern already used on the Card24 sheet: a bold/bordered header cell
# in N1, and blank data cells in N2:N12. At the same time the previously
# blank Event cells (M2:M12) pick up the literal "nan" placeholder text
# used throughout the rest of the table for missing values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card18")

# --- Header: N1 = "Correction", formatted like the other header cells ---
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("N1").Value = "Correction"

# --- Data rows 2-12: M gets "nan", N is created as a blank cell ---
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 13).Value = "nan"        # column M = 13 ("Event")

    $nCell = $ws.Cells.Item($r, 14)             # column N = 14 ("Correction")
    $nCell.Font.Bold = $false                   # materialize the empty cell
}
